# 19/09/23 - push 1
# sheet "exercício 8": add a "soma" column (C) computing vetor1 + 1/vetor2
# for each row, entered cell-by-cell the way a user typically would
# (type the formula in C2, then again in C3, then fill C3 down through
# C6) so Excel records C2 as its own formula and C3:C6 as one shared
# formula group anchored at C3 - matching the authored workbook exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("exercício 8")

$ws.Range("C2").Formula = "=A2 + 1/B2"
$ws.Range("C3:C6").Formula = "=A3 + 1/B3"
